$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Terms fixed for customers" - the terms column (AB) held a plain numeric
# placeholder (1 / 2); replace it with the real terms-of-payment text for
# each customer record.
$ws.Range("AB1").Value = "terms1"
$ws.Range("AB2").Value = "terms2"

# Scroll/select like the author left the sheet: top-left visible cell at
# AA1, with AB3 as the active selection.
$excel.Goto($ws.Range("AA1"), $false)
$ws.Range("AB3").Select()
